$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (Las Pistol): D2 40/B2 -> 50/B2, F2 35/B2 -> 40/B2
$ws.Range("D2").Formula = "=50/B2"
$ws.Range("F2").Formula = "=40/B2"

# Row 4 (Machine Las Pistol): D4 45/B4 -> 50/B4, F4 35/B4 -> 40/B4
$ws.Range("D4").Formula = "=50/B4"
$ws.Range("F4").Formula = "=40/B4"

# Update the active selection to mirror the saved view state (F5)
$ws.Range("F5").Select()
